# Adds "Call Broadcast" and "get Recorded audio file" API rows (row 36 and
# row 39) to the Auth / Voice-bot API sheet, matching the author's commit:
# "Added Call Broadcast API and get Recorded audio files in Auth
# application for voice bot"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Exact text content (kept as single-quoted here-strings so nothing is
# expanded/escaped by the shell).
# ---------------------------------------------------------------------

$urlCallBroadcast = @'
http://1msg.1point1.in:3001/api/auth/j-v1/call/Broadcast/
'@

$payloadCallBroadcast = @'
{
    "agent_id": 1,
    "user_id":1,
    "group_id":1,
    "phone_numbers": [
      { 
          "name": "naga", "number": "9821209237","contact_id":"1" 
           }
     ]
  }
'@

$curlCallBroadcast = @'
curl --location 'http://1msg.1point1.in:3001/api/auth/j-v1/call/Broadcast/' \
--header 'Content-Type: application/json' \
--data '{
    "agent_id": 1,
    "user_id":1,
    "group_id":1,
    "phone_numbers": [
      { 
          "name": "naga", "number": "9821209237","contact_id":"1" 
           }
     ]
  }
 '
'@

$urlGetRecordedAudio = @'
http://1msg.1point1.in:3001/api/auth/j-v1/get/recorded/audio/file/
'@

$titleGetRecordedAudio = @'
get Recorded audio file     ( POST ) 
'@

$titleCallBroadcast = @'
Call Broadcast   ( POST )
'@

$payloadGetRecordedAudio = @'
{
    "file_path": "/usr/share/freeswitch/var/lib/freeswitch/recordings/c6b69a41-dec9-4257-b061-50b3da35b3e6.wav"
}
'@

$curlGetRecordedAudio = @'
curl --location 'http://1msg.1point1.in:3001/api/auth/j-v1/get/recorded/audio/file/' \
--header 'Content-Type: application/json' \
--data '{
    "file_path": "/usr/share/freeswitch/var/lib/freeswitch/recordings/c6b69a41-dec9-4257-b061-50b3da35b3e6.wav"
}'
'@

# ---------------------------------------------------------------------
# Row 36 - "Call Broadcast ( POST )"
# ---------------------------------------------------------------------

# B36: URL + hyperlink, formatted like the other hyperlinked URL cells
# (C2 / C17 use style index 5 - Hyperlink + wrap + vertical-center).
$ws.Range("B36").Value = $urlCallBroadcast
$ws.Hyperlinks.Add($ws.Range("B36"), $urlCallBroadcast)
$ws.Range("C2").Copy()
$ws.Range("B36").PasteSpecial(-4122)

# C36 / D36: payload + curl command, formatted like the other payload /
# response cells (D28 uses style index 6 - wrap only).
$ws.Range("D28").Copy()
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C36").Value = $payloadCallBroadcast

$ws.Range("D28").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D36").Value = $curlCallBroadcast

# ---------------------------------------------------------------------
# Row 39 - "get Recorded audio file ( POST )"
# ---------------------------------------------------------------------

# B39: URL + hyperlink (new, unformatted cell -> Excel creates a fresh
# Hyperlink-only style for it).
$ws.Range("B39").Value = $urlGetRecordedAudio
$ws.Hyperlinks.Add($ws.Range("B39"), $urlGetRecordedAudio)

# A39: plain title text, default styling (same as a never-formatted cell).
$ws.Range("A39").Value = $titleGetRecordedAudio

# A36: title text, formatted like A2 (wrap + vertical-center).
$ws.Range("A2").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value = $titleCallBroadcast

# C39 / D39: payload + curl command, same formatting family as C36/D36.
$ws.Range("D28").Copy()
$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("C39").Value = $payloadGetRecordedAudio

$ws.Range("D28").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("D39").Value = $curlGetRecordedAudio

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Row heights for the two new rows.
# ---------------------------------------------------------------------

$ws.Rows.Item(36).RowHeight = 248.4
$ws.Rows.Item(39).RowHeight = 110.4

# ---------------------------------------------------------------------
# Selection, matching where the author ended up.
# ---------------------------------------------------------------------

$ws.Range("D39").Select()
